# Database_Thresholds.xlsx — refresh Quantiles with 7-15-24 exports
# (WQ_Discrete re-exported 7/11/24)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Updated: 2024-07-10" -> "Updated: 2024-07-15" ---
$ws.Range("A3").Value = "Updated: 2024-07-15"

# --- Single recalculated quantile value ---
$ws.Range("O68").Value = 3.139395

$oldCommitString = "Database_Thresholds.xlsx, Git Commit ID: ea27439d179b540cda00ced4b43f858e55a0f99d"
$newCommitString = "Database_Thresholds.xlsx, Git Commit ID: 51214cf344547100c5a54d26465d79cec9558701"

$oldDate = 45483
$newDate = 45488

# --- Data rows 8-96: ScriptLatestRunDate (X) refreshed for every row ---
for ($r = 8; $r -le 96; $r++) {
    $xCell = $ws.Cells.Item($r, 24)
    if ($xCell.Value2 -eq $oldDate) {
        $xCell.Value = $newDate
    }

    # ActionNeededDate (T) and QuantileDate (Y) only where populated
    $tCell = $ws.Cells.Item($r, 20)
    if ($tCell.Value2 -eq $oldDate) {
        $tCell.Value = $newDate
    }

    $yCell = $ws.Cells.Item($r, 25)
    if ($yCell.Value2 -eq $oldDate) {
        $yCell.Value = $newDate
    }

    # QuantileSource (U) commit-id string, only where it references the old commit
    $uCell = $ws.Cells.Item($r, 21)
    if ($uCell.Value2 -eq $oldCommitString) {
        $uCell.Value = $newCommitString
    }
}
